$wb = $excel.ActiveWorkbook

# ----- Sheet ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 11000
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H62").Value = 3301.6667
$ws.Range("I62").Value = 3452.5
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 3452.5
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -2828.5
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 3301.6667
$ws.Range("I65").Value = 3452.5
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 17262.5
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -14142.5
$ws.Range("N65").Value = -21240
$ws.Range("H98").Value = 652.2963
$ws.Range("I98").Value = 604.3200000000001
$ws.Range("K98").Value = 604.3200000000001
$ws.Range("M98").Value = 893.6799999999999
$ws.Range("H107").Value = 1008.7273
$ws.Range("I107").Value = 1162.9333
$ws.Range("J107").Value = 678.2857
$ws.Range("K107").Value = 1162.9333
$ws.Range("L107").Value = 678.2857
$ws.Range("M107").Value = 757.0667000000001
$ws.Range("N107").Value = -4518.2857
$ws.Range("H116").Value = 3837.4736
$ws.Range("I116").Value = 2012.5
$ws.Range("K116").Value = 2012.5
$ws.Range("M116").Value = 1429.5
$ws.Range("H122").Value = 652.2963
$ws.Range("I122").Value = 604.3200000000001
$ws.Range("K122").Value = 1812.96
$ws.Range("M122").Value = 637.04
$ws.Range("H125").Value = 1361.4286
$ws.Range("I125").Value = 400
$ws.Range("J125").Value = 2082.5
$ws.Range("K125").Value = 3600
$ws.Range("L125").Value = 18742.5
$ws.Range("M125").Value = -1140
$ws.Range("N125").Value = -23662.5
$ws.Range("H132").Value = 2197.838
$ws.Range("I132").Value = 2303.5293
$ws.Range("K132").Value = 6910.5879
$ws.Range("M132").Value = -4380.5879

# ----- Sheet ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3747.8
$ws.Range("I102").Value = 1999.6666
$ws.Range("J102").Value = 6370
$ws.Range("K102").Value = 1999.6666
$ws.Range("L102").Value = 6370
$ws.Range("M102").Value = -377.6666
$ws.Range("N102").Value = -9614
$ws.Range("H107").Value = 49999.668
$ws.Range("J107").Value = 49999.668
$ws.Range("L107").Value = 49999.668
$ws.Range("N107").Value = -57679.668
$ws.Range("H122").Value = 1959.2727
$ws.Range("I122").Value = 1950.3334
$ws.Range("K122").Value = 5851.0002
$ws.Range("M122").Value = -3401.0002

# ----- Sheet BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 939.82355
$ws.Range("I20").Value = 1188.5264
$ws.Range("J20").Value = 624.8
$ws.Range("K20").Value = 1188.5264
$ws.Range("L20").Value = 624.8
$ws.Range("M20").Value = -941.5264
$ws.Range("N20").Value = -1118.8
$ws.Range("H80").Value = 727.7273
$ws.Range("I80").Value = 719.6
$ws.Range("J80").Value = 734.5
$ws.Range("K80").Value = 719.6
$ws.Range("L80").Value = 734.5
$ws.Range("M80").Value = 278.4
$ws.Range("N80").Value = -2730.5
$ws.Range("H83").Value = 727.7273
$ws.Range("I83").Value = 719.6
$ws.Range("J83").Value = 734.5
$ws.Range("K83").Value = 3598
$ws.Range("L83").Value = 3672.5
$ws.Range("M83").Value = 1394
$ws.Range("N83").Value = -13656.5
$ws.Range("H140").Value = 41396.57
$ws.Range("J140").Value = 41396.57
$ws.Range("L140").Value = 41396.57
$ws.Range("N140").Value = -51756.57

# ----- Sheet CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12340.324
$ws.Range("I31").Value = 15711.538
$ws.Range("J31").Value = 4372
$ws.Range("K31").Value = 15711.538
$ws.Range("L31").Value = 4372
$ws.Range("M31").Value = -15416.538
$ws.Range("N31").Value = -4962
$ws.Range("H34").Value = 12340.324
$ws.Range("I34").Value = 15711.538
$ws.Range("J34").Value = 4372
$ws.Range("K34").Value = 15711.538
$ws.Range("L34").Value = 4372
$ws.Range("M34").Value = -15509.538
$ws.Range("N34").Value = -4776
$ws.Range("H99").Value = 4103.24
$ws.Range("J99").Value = 5914.875
$ws.Range("L99").Value = 5914.875
$ws.Range("N99").Value = -8910.875
$ws.Range("H122").Value = 1122.5
$ws.Range("J122").Value = 1016.75
$ws.Range("L122").Value = 3050.25
$ws.Range("N122").Value = -7950.25
$ws.Range("H126").Value = 4103.24
$ws.Range("J126").Value = 5914.875
$ws.Range("L126").Value = 17744.625
$ws.Range("N126").Value = -22684.625
$ws.Range("H134").Value = 1255.0834
$ws.Range("I134").Value = 1007.75
$ws.Range("K134").Value = 3023.25
$ws.Range("M134").Value = -488.25
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# ----- Sheet CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H131").Value = 781.0599999999999
$ws.Range("I131").Value = 700
$ws.Range("J131").Value = 781.8788
$ws.Range("K131").Value = 2100
$ws.Range("L131").Value = 2345.6364
$ws.Range("M131").Value = 2940
$ws.Range("N131").Value = -12425.6364

# ----- Sheet GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 94.22221999999999
$ws.Range("I2").Value = 60.42857
$ws.Range("J2").Value = 212.5
$ws.Range("K2").Value = 60.42857
$ws.Range("L2").Value = 212.5
$ws.Range("M2").Value = 52.57143
$ws.Range("N2").Value = -438.5
$ws.Range("H43").Value = 2197.4
$ws.Range("I43").Value = 2197.4
$ws.Range("K43").Value = 2197.4
$ws.Range("M43").Value = -2046.4
$ws.Range("H46").Value = 5270.5
$ws.Range("I46").Value = 5270.5
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 5270.5
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -5114.5
$ws.Range("N46").ClearContents()
$ws.Range("H57").Value = 29445
$ws.Range("J57").Value = 29445
$ws.Range("L57").Value = 29445
$ws.Range("N57").Value = -31085
$ws.Range("H80").Value = 4052.2222
$ws.Range("I80").Value = 3663.3333
$ws.Range("J80").Value = 4246.6665
$ws.Range("K80").Value = 3663.3333
$ws.Range("L80").Value = 4246.6665
$ws.Range("M80").Value = -2665.3333
$ws.Range("N80").Value = -6242.6665
$ws.Range("H83").Value = 4052.2222
$ws.Range("I83").Value = 3663.3333
$ws.Range("J83").Value = 4246.6665
$ws.Range("K83").Value = 18316.6665
$ws.Range("L83").Value = 21233.3325
$ws.Range("M83").Value = -13324.6665
$ws.Range("N83").Value = -31217.3325
$ws.Range("H132").Value = 154641.2
$ws.Range("I132").Value = 130301.625
$ws.Range("J132").Value = 251999.5
$ws.Range("K132").Value = 390904.875
$ws.Range("L132").Value = 755998.5
$ws.Range("M132").Value = -388374.875
$ws.Range("N132").Value = -761058.5
$ws.Range("H140").Value = 52500
$ws.Range("J140").Value = 52500
$ws.Range("L140").Value = 52500
$ws.Range("N140").Value = -62860

# ----- Sheet LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4069.8
$ws.Range("I7").Value = 3729
$ws.Range("J7").Value = 5433
$ws.Range("K7").Value = 3729
$ws.Range("L7").Value = 5433
$ws.Range("M7").Value = -3617
$ws.Range("N7").Value = -5657
$ws.Range("H16").Value = 458.41666
$ws.Range("I16").Value = 431.375
$ws.Range("J16").Value = 512.5
$ws.Range("K16").Value = 431.375
$ws.Range("L16").Value = 512.5
$ws.Range("M16").Value = -261.375
$ws.Range("N16").Value = -852.5
$ws.Range("H55").Value = 312.66666
$ws.Range("I55").Value = 142
$ws.Range("J55").Value = 526
$ws.Range("K55").Value = 142
$ws.Range("L55").Value = 526
$ws.Range("M55").Value = 31
$ws.Range("N55").Value = -872
$ws.Range("H93").Value = 2634.0527
$ws.Range("I93").Value = 2952.5386
$ws.Range("J93").Value = 1944
$ws.Range("K93").Value = 2952.5386
$ws.Range("L93").Value = 1944
$ws.Range("M93").Value = -1704.5386
$ws.Range("N93").Value = -4440
$ws.Range("H126").Value = 4069.8
$ws.Range("I126").Value = 3729
$ws.Range("J126").Value = 5433
$ws.Range("K126").Value = 11187
$ws.Range("L126").Value = 16299
$ws.Range("M126").Value = -8717
$ws.Range("N126").Value = -21239
$ws.Range("H132").Value = 862917
$ws.Range("I132").Value = 929133.7
$ws.Range("K132").Value = 2787401.1
$ws.Range("M132").Value = -2784871.1

# ----- Sheet WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 39929
$ws.Range("J123").Value = 39929
$ws.Range("L123").Value = 39929
